# Hotel_Management_System RTM_Group_2.xlsx - "Add files via upload" edit
# Updates the Req / Design Mapping columns (B, C) and populates the new
# UT Mapping / IT Mapping columns (E, F) for rows 2-7, then moves the
# active selection to F7 (matching the author's final cursor position).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Rachala Manisha / admin_validation()
$ws.Cells.Item(2, 2).Value = "[HMS]//09/10"
$ws.Cells.Item(2, 3).Value = "2.1.8"
$ws.Cells.Item(2, 5).Value = "ADMIN_VALIDATION_007"
$ws.Cells.Item(2, 6).Value = "admin_validation_007"

# Row 3 - Aishwarya Palled / add_room(),display_room(),exit()
$ws.Cells.Item(3, 2).Value = "[HMS]/02/04"
$ws.Cells.Item(3, 3).Value = "2.1.2,2.1.3,2.1.7"
$ws.Cells.Item(3, 5).Value = "ADD_ROOM_001,DISPLAY_ROOM_002,EXIT_006"
$ws.Cells.Item(3, 6).Value = "add_room_001,display_room_002,exit_006"

# Row 4 - Swathi B R / modify_room(),search_customer(),view_customer()
$ws.Cells.Item(4, 2).Value = "[HMS]/07"
$ws.Cells.Item(4, 3).Value = "2.1.4,2.1.5,2.1.6"
$ws.Cells.Item(4, 5).Value = "MODIFY_ROOM_003,SEARCH_CUSTOMER_004,VIEW_CUSTOMER_005"
$ws.Cells.Item(4, 6).Value = "modify_room_003,search_customer_004,view_customer_005"

# Row 5 - Shridevi Prabhu S / register_new_customer(),book_room()
$ws.Cells.Item(5, 2).Value = "[HMS]/03/11"
$ws.Cells.Item(5, 3).Value = "2.2.1,2.2.2"
$ws.Cells.Item(5, 5).Value = "REGISTER_NEW_CUSTOMER_008,BOOK_ROOM_008"
$ws.Cells.Item(5, 6).Value = "register_new_customer_008,book_room_009"

# Row 6 - Geethanjali Goddumarri / search_room(),checkout_room,exit()
$ws.Cells.Item(6, 2).Value = "[HMS]/04/11"
$ws.Cells.Item(6, 3).Value = "2.2.3,2.2.4,2.2.5"
$ws.Cells.Item(6, 5).Value = "SEARCH_ROOM_010,CHECKOUT_ROOM_011,EXIT_012"
$ws.Cells.Item(6, 6).Value = "search_room_010,checkout_room_011,exit_13"

# Row 7 - Mallela Reddypriya / admin_login()
$ws.Cells.Item(7, 2).Value = "[HMS]/01/05"
$ws.Cells.Item(7, 3).Value = "2.1.1"
$ws.Cells.Item(7, 5).Value = "ADMIN_VALIDATION_013"
$ws.Cells.Item(7, 6).Value = "admin_login_14"

# Move the selection to match the saved cursor position
[void]$ws.Range("F7").Select()
